$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# "Stored as a binary string of length 11 instead of 11 columns of data"
#   -> "Stored as a binary string of length 12 instead of 12 columns of data"
$para2 = $tr.Paragraphs(2, 1)
$para2.Runs(1).Text = "Stored as a binary string of length 12 instead of 12 columns of data"

# "Value of ‘1’ indicates the flag is set.  e.g. 00100000000"
#   -> "Value of ‘1’ indicates the flag is set.  e.g. 001000000000"
$para3 = $tr.Paragraphs(3, 1)
$para3.Runs(1).Text = "Value of ‘1’ indicates the flag is set.  e.g. 001000000000"

# "All combinations can be represented as a number from 1 to 2048 (i.e. 2" + "11" (superscript) + "-1)..."
#   -> "All combinations can be represented as a number from 0 to 4095 (i.e. 2" + "12" (superscript) + "-1)..."
$para4 = $tr.Paragraphs(4, 1)
$para4.Runs(1).Text = "All combinations can be represented as a number from 0 to 4095 (i.e. 2"
$para4.Runs(2).Text = "12"
